$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.096.04'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.652.18'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.39'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5203'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2644'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06335'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.40'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07687'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.594'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.653.28'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.878.85'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5594'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8149'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.40'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.106.67'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.15%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.629'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.48'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '191.69'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.923'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.003'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.18'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.60%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.219'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.71%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.504'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05488'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.268'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.442'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.359'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.422'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9485'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.787'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5638'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01582'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.850'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.64%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.028.60'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8282'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.45%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.28'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.794.38'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.55'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.67%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₈109'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +4.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9993'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.001'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4336'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05166'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.69%  '
